# Fix several sentences with incorrect German/Polish rektion examples.
# (Matches commit: "Poprawione zdania z rekcja" - corrected sentences with rektion)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arkusz1")

$ws.Range("A38").Value  = "Die Anzahl der Teilnehmer wurde auf zwanzig begrenzt."
$ws.Range("A75").Value  = "Wir haben uns für diesen Termin entschieden."
$ws.Range("B122").Value = "Bardzo się zdenerwowaliśmy tym incydentem. (perf.)"
$ws.Range("A139").Value = "Er sollte sich darüber nicht wundern."
$ws.Range("B157").Value = "Kiedy się o tym dowiedziałeś? (perf.)"
$ws.Range("B161").Value = "Nic nie słyszałem o tej książce. (perf.)"
$ws.Range("B163").Value = "Dużo skorzystałem na tej praktyce. (perf.)"
$ws.Range("B174").Value = "Odwiedziny przyczyniły się do pogłębienia relacji. (perf.)"
$ws.Range("B179").Value = "Ten kraj rozwinął się w kraj przemysłowy. (perf.)"
$ws.Range("B180").Value = "Do czego doprowadziła ta polityka? (perf.)"
$ws.Range("B195").Value = "Niektórzy ludzie brzydzą się pająka."
$ws.Range("B200").Value = "Peter rozstał się ze swoim przyjacielem. (perf.)"
$ws.Range("B202").Value = "Dziś Jens zaręczył się ze swoją dziewczyną. (perf.)"
$ws.Range("A205").Value = "Die Qualität des neuen Artikels weicht von der des alten ab."
$ws.Range("A206").Value = "Ich bin neugierig daran, was er sagen wird."
$ws.Range("A224").Value = "Sie staunen über die Architektur der Stadt."
$ws.Range("B226").Value = "On został wykluczony z partii. (perf.)"
$ws.Range("A236").Value = "Sie wird sich nicht dazu äußern."
$ws.Range("A237").Value = "Sie hat sich sehr für die Rechte von Kindern engagiert."
$ws.Range("B238").Value = "Ona przestraszyła się własnego cienia. (perf.)"
$ws.Range("B248").Value = "Ukryła się (swoją twarz) za gazetą. (perf.)"
$ws.Range("B252").Value = "Każdy narzeka z powodu upału."

# Move the view/selection to where the user ended editing.
$ws.Application.ActiveWindow.ScrollRow = 113
$ws.Range("B123").Select()
